$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 28, shifting the existing rows 28-50 down to 30-52.
$ws.Rows("28:29").Insert()

# New "Sandia" records (matches the commit's weekly price update).
$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(28, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(28, 4).Value = 44904
$ws.Cells.Item(28, 5).Value = 15
$ws.Cells.Item(28, 6).Value = 100112028
$ws.Cells.Item(28, 7).Value = "Sandia"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Segunda"
$ws.Cells.Item(28, 10).Value = 900
$ws.Cells.Item(28, 11).Value = 470
$ws.Cells.Item(28, 12).Value = 490
$ws.Cells.Item(28, 13).Value = 481
$ws.Cells.Item(28, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(28, 15).Value = "Perú"
$ws.Cells.Item(28, 16).Value = 481
$ws.Cells.Item(28, 17).Value = 1
$ws.Cells.Item(28, 18).Value = "Hortaliza"

$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(29, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(29, 4).Value = 44904
$ws.Cells.Item(29, 5).Value = 15
$ws.Cells.Item(29, 6).Value = 100112028
$ws.Cells.Item(29, 7).Value = "Sandia"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Tercera"
$ws.Cells.Item(29, 10).Value = 700
$ws.Cells.Item(29, 11).Value = 470
$ws.Cells.Item(29, 12).Value = 490
$ws.Cells.Item(29, 13).Value = 479
$ws.Cells.Item(29, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(29, 15).Value = "Perú"
$ws.Cells.Item(29, 16).Value = 479
$ws.Cells.Item(29, 17).Value = 1
$ws.Cells.Item(29, 18).Value = "Hortaliza"
